$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(78.7,  0.01251, 0.00625, 1.39921)
    3  = @(63,    0.04376, 0.02188, 3.42212)
    4  = @(41.3,  0.06219, 0.01555, 3.59001)
    5  = @(43.75, 0.02716, 0.00679, 1.85701)
    6  = @(24.3,  0.03712, 0.00619, 1.74744)
    7  = @(28.95, 0.07126, 0.01188, 3.25552)
    8  = @(21.2,  0.09186, 0.01148, 3.09895)
    9  = @(20,    0.04217, 0.00527, 1.74764)
    10 = @(13.6,  0.03833, 0.00383, 1.25236)
    11 = @(14.55, 0.11049, 0.01105, 2.82535)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 3).Value = $rowValues[0]
    $ws.Cells.Item($row, 4).Value = $rowValues[1]
    $ws.Cells.Item($row, 5).Value = $rowValues[2]
    $ws.Cells.Item($row, 6).Value = $rowValues[3]
}
